$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per commit diff
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.107.60'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.655.06'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.08'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5256'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2609'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06354'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.43'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07806'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.506'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.649.84'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5483'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0₅8235'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.44'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.05%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.122.44'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.580'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '191.67'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.038'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '141.89'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1247'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.265'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.430'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.05919'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.279'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.525'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.258'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.590'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.791'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.411'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.5706'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01620'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.8505'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.795'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.31%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '103.07'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.029.94'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.801.00'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '57.22'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.9992'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4299'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.11%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.478'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05166'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.840'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.09720'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.59%  '
